# Update "Pais" (countries) data sheet with newer COVID-19 stats snapshot
# and bump the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados..." timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 18:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1711008
$ws.Range("C4").Value = 4782
$ws.Range("D4").Value = 467134
$ws.Range("E4").Value = 1143924
$ws.Range("G4").Value = 145
$ws.Range("H4").Value = 99950

# Row 5 - Brasil
$ws.Range("B5").Value = 377711
$ws.Range("C5").Value = 1042
$ws.Range("E5").Value = 200272
$ws.Range("G5").Value = 84
$ws.Range("H5").Value = 23606

# Row 9 - Italia
$ws.Range("B9").Value = 230555
$ws.Range("C9").Value = 397
$ws.Range("D9").Value = 144658
$ws.Range("E9").Value = 52942
$ws.Range("G9").Value = 78
$ws.Range("H9").Value = 32955

# Row 13 - India
$ws.Range("B13").Value = 150313
$ws.Range("C13").Value = 5363
$ws.Range("D13").Value = 63536
$ws.Range("E13").Value = 82443
$ws.Range("G13").Value = 162
$ws.Range("H13").Value = 4334

# Row 69 - Irak
$ws.Range("B69").Value = 4848
$ws.Range("C69").Value = 216
$ws.Range("E69").Value = 1868
$ws.Range("G69").Value = 6
$ws.Range("H69").Value = 169

# Row 70 - Azerbaiyan
$ws.Range("B70").Value = 4403
$ws.Range("C70").Value = 132
$ws.Range("D70").Value = 2819
$ws.Range("E70").Value = 1532
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 52

# Row 72 - Luxemburgo
$ws.Range("B72").Value = 3995
$ws.Range("C72").Value = 2
$ws.Range("D72").Value = 3783

# Row 96 - Mayotte
$ws.Range("B96").Value = 1634
$ws.Range("C96").Value = 25
$ws.Range("E96").Value = 720

# Row 140 - Cabo Verde
$ws.Range("E140").Value = 231
$ws.Range("G140").Value = 1
$ws.Range("H140").Value = 4
